$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Partial Molar Volume Analysis")
$ws3 = $wb.Worksheets.Item("Also with uncertainty!")

# --- Sheet "Partial Molar Volume Analysis": refreshed polynomial-fit results ---
$ws2.Range("B2").Value = 0.0
$ws2.Range("F2").Value = 0.01053987096240326
$ws2.Range("H2").Value = 51.21948003498717
$ws2.Range("F3").Value = 0.011113095737457345
$ws2.Range("H3").Value = 54.00511896986348
$ws2.Range("B4").Value = 4.224784285543178
$ws2.Range("F4").Value = 0.011370306108005027
$ws2.Range("H4").Value = 55.25505661009192
$ws2.Range("B5").Value = 5.981995448556714
$ws2.Range("F5").Value = 0.011461576737019237
$ws2.Range("H5").Value = 55.698594692982084
$ws2.Range("B6").Value = 7.944837705114385
$ws2.Range("F6").Value = 0.011607112618723691
$ws2.Range("H6").Value = 56.405839801952176
$ws2.Range("B7").Value = 9.888986225895316
$ws2.Range("F7").Value = 0.011959847739864717
$ws2.Range("H7").Value = 58.11998882326069
$ws2.Range("F8").Value = 0.01264874296840362
$ws2.Range("H8").Value = 61.467739050014956

# --- Sheet "Also with uncertainty!": refreshed polynomial-fit results ---
$ws3.Range("B2").Value = 0.0
$ws3.Range("C2").Value = 0.004673433944184932
$ws3.Range("E2").Value = 0.0023503417324359246
$ws3.Range("F2").Value = 0.010539870962403301
$ws3.Range("G2").Value = 0.0013482065001652753
$ws3.Range("H2").Value = 51.21948003498738
$ws3.Range("I2").Value = 6.5752203803287195
$ws3.Range("C3").Value = 0.022478706378621635
$ws3.Range("E3").Value = 0.0023503417324359246
$ws3.Range("F3").Value = 0.011113095737457345
$ws3.Range("G3").Value = 0.00030566353984542554
$ws3.Range("H3").Value = 54.00511896986348
$ws3.Range("I3").Value = 1.596608553343804
$ws3.Range("B4").Value = 4.224784285543178
$ws3.Range("C4").Value = 0.0460367269257998
$ws3.Range("E4").Value = 0.0023503417324359246
$ws3.Range("F4").Value = 0.011370306108005017
$ws3.Range("G4").Value = 0.00036816970517249717
$ws3.Range("H4").Value = 55.25505661009188
$ws3.Range("I4").Value = 1.8867614750867665
$ws3.Range("B5").Value = 5.981995448556714
$ws3.Range("C5").Value = 0.06501618328409213
$ws3.Range("E5").Value = 0.0023503417324359246
$ws3.Range("F5").Value = 0.011461576737019233
$ws3.Range("G5").Value = 0.0002461288425864187
$ws3.Range("H5").Value = 55.69859469298207
$ws3.Range("I5").Value = 1.3398502660856424
$ws3.Range("B6").Value = 7.944837705114385
$ws3.Range("C6").Value = 0.08625295313143291
$ws3.Range("E6").Value = 0.0023503417324359246
$ws3.Range("F6").Value = 0.011607112618723695
$ws3.Range("G6").Value = 0.0003781153869899814
$ws3.Range("H6").Value = 56.4058398019522
$ws3.Range("I6").Value = 1.9365567949520686
$ws3.Range("B7").Value = 9.888986225895316
$ws3.Range("E7").Value = 0.0023503417324359246
$ws3.Range("F7").Value = 0.011959847739864716
$ws3.Range("G7").Value = 0.0003080436305770664
$ws3.Range("H7").Value = 58.119988823260684
$ws3.Range("I7").Value = 1.6241532556642557
$ws3.Range("E8").Value = 0.0023503417324359246
$ws3.Range("F8").Value = 0.012648742968403584
$ws3.Range("G8").Value = 0.0013179241631915053
$ws3.Range("H8").Value = 61.467739050014785
$ws3.Range("I8").Value = 6.439145069972003

# --- Active sheet moves back to Sheet1 (matches updated tabSelected/activeTab state) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
